# Unit correction: the "dia" (column D) and "end_diameter" (column E)
# measurements were recorded one order of magnitude too large. The
# author's commit message ("Re-upload due to uinit correction") and the
# data diff show every value in D2:E121 divided by 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 121; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dCell.Value = $dCell.Value() / 10
    $eCell.Value = $eCell.Value() / 10
}

# Carry over the author's final selection/scroll state from the re-upload.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("K9").Select() | Out-Null
